# Updated cryptos list on Sat May 27 07:44:15 UTC 2023 with GitHub Actions
#
# Refreshes the Price / Volume(1h) snapshot for every coin row, and
# corrects rows 11-13 whose Coin/Link/Price/Volume content had shifted
# out of rank order (Solana/TRON/WrappedEther -> WrappedEther/Solana/TRON).
#
# Every new value is written with a leading apostrophe so Excel treats
# numeric-looking strings (e.g. '309.51', '66.00') as literal text,
# matching the workbook's original inline-string / General-format cells
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.943.13"
$ws.Range("E2").Value = "'  +1.23%  "
$ws.Range("D3").Value = "'1.846.85"
$ws.Range("E3").Value = "'  +1.27%  "
$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("D5").Value = "'309.51"
$ws.Range("E5").Value = "'  +0.54%  "
$ws.Range("E6").Value = "'  +0.28%  "
$ws.Range("D7").Value = "'0.4775"
$ws.Range("E7").Value = "'  +2.84%  "
$ws.Range("E8").Value = "'  +1.83%  "
$ws.Range("D9").Value = "'0.07220"
$ws.Range("E9").Value = "'  +1.28%  "
$ws.Range("D10").Value = "'0.9277"
$ws.Range("B11").Value = "'WrappedEther"
$ws.Range("C11").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "'1.971.44"
$ws.Range("E11").Value = "'  +8.05%  "
$ws.Range("B12").Value = "'Solana"
$ws.Range("C12").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'19.73"
$ws.Range("E12").Value = "'  +1.93%  "
$ws.Range("B13").Value = "'TRON"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07693"
$ws.Range("E13").Value = "'  -0.84%  "
$ws.Range("D14").Value = "'5.319"
$ws.Range("E14").Value = "'  +1.08%  "
$ws.Range("D15").Value = "'6.409"
$ws.Range("E15").Value = "'  +1.18%  "
$ws.Range("D16").Value = "'88.93"
$ws.Range("E16").Value = "'  +1.43%  "
$ws.Range("E17").Value = "'  +0.27%  "
$ws.Range("D18").Value = "'0.000008640"
$ws.Range("E18").Value = "'  +1.00%  "
$ws.Range("E19").Value = "'  +0.28%  "
$ws.Range("D20").Value = "'26.982.01"
$ws.Range("E20").Value = "'  +1.24%  "
$ws.Range("D21").Value = "'14.56"
$ws.Range("E21").Value = "'  +2.83%  "
$ws.Range("D22").Value = "'5.056"
$ws.Range("E22").Value = "'  +0.85%  "
$ws.Range("D23").Value = "'10.64"
$ws.Range("E23").Value = "'  +1.10%  "
$ws.Range("D24").Value = "'1.923"
$ws.Range("E24").Value = "'  +0.28%  "
$ws.Range("D25").Value = "'152.61"
$ws.Range("E25").Value = "'  +0.18%  "
$ws.Range("E26").Value = "'  +1.56%  "
$ws.Range("D27").Value = "'1.999"
$ws.Range("E27").Value = "'  +1.77%  "
$ws.Range("D28").Value = "'114.30"
$ws.Range("E28").Value = "'  +0.50%  "
$ws.Range("D29").Value = "'4.923"
$ws.Range("E29").Value = "'  +1.64%  "
$ws.Range("D30").Value = "'0.08877"
$ws.Range("E30").Value = "'  +0.84%  "
$ws.Range("D31").Value = "'3.312"
$ws.Range("E31").Value = "'  +5.63%  "
$ws.Range("E32").Value = "'  +3.40%  "
$ws.Range("D33").Value = "'0.7466"
$ws.Range("E33").Value = "'  +1.95%  "
$ws.Range("D34").Value = "'4.492"
$ws.Range("E34").Value = "'  +1.30%  "
$ws.Range("D35").Value = "'2.726"
$ws.Range("E35").Value = "'  +0.13%  "
$ws.Range("D36").Value = "'1.108"
$ws.Range("E36").Value = "'  +3.05%  "
$ws.Range("D37").Value = "'0.01959"
$ws.Range("E37").Value = "'  +2.06%  "
$ws.Range("D38").Value = "'0.05262"
$ws.Range("E38").Value = "'  +2.72%  "
$ws.Range("E39").Value = "'  +1.59%  "
$ws.Range("D40").Value = "'0.5197"
$ws.Range("E40").Value = "'  +2.83%  "
$ws.Range("D41").Value = "'6.963"
$ws.Range("E41").Value = "'  +1.36%  "
$ws.Range("D42").Value = "'0.1509"
$ws.Range("E42").Value = "'  +1.09%  "
$ws.Range("D43").Value = "'8.202"
$ws.Range("E43").Value = "'  +2.66%  "
$ws.Range("D44").Value = "'10.54"
$ws.Range("E44").Value = "'  +4.91%  "
$ws.Range("D45").Value = "'0.4725"
$ws.Range("E45").Value = "'  +1.55%  "
$ws.Range("E46").Value = "'  +0.36%  "
$ws.Range("D47").Value = "'101.46"
$ws.Range("E47").Value = "'  +3.18%  "
$ws.Range("E48").Value = "'  +2.95%  "
$ws.Range("D49").Value = "'66.00"
$ws.Range("E49").Value = "'  +3.68%  "
$ws.Range("D50").Value = "'0.06017"
$ws.Range("E50").Value = "'  +0.55%  "
$ws.Range("D51").Value = "'0.8855"
$ws.Range("E51").Value = "'  +3.94%  "
